$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1235
$ws1.Range("F6").Value = 156

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1235
$ws4.Range("F6").Value = 156
